# Extend the table from columns A:O to A:Q (add two new columns, P and Q)
# and flip a few of the "contingency" indicator columns (I, K, M, O) for
# every data row, per the "contingencies with rene fine" change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells P1/Q1 continue the 0..13 sequence already in B1:O1 and
# should carry the same bold/centered/bordered header style as O1.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Data rows 2-25: swap values in columns I/K and M/O, and append two new
# columns (P, Q) both populated with 2.
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
